# Handback report generation:
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet (zh-cn / de-de columns) and on each language detail
#   sheet's Status column.
# - Each language detail sheet now carries the handback target file, the
#   handback xliff file name, and the handback datetime, plus a hyperlink on
#   the newly-populated "Latest Target File" cell.
# - A couple of columns are widened to fit the longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f77bbdd5ba31ca8bacd1e92c013abfabe4efea39/e2e"

# --- Overview sheet: refresh the per-language status summary -------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de columns to fit the new, longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

# --- zh-cn detail sheet ----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-30 16:44:49"
$wsZh.Range("K3").Value = "2016-08-30 16:44:49"

# Re-create the hyperlinks in document order (A2, I2, A3, I3) so the new
# "Latest Target File" links for row 2 / row 3 sit next to the existing
# "Source File Name" links.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBlobBase/a.md", "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$repoBlobBase/a.md", "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBlobBase/b.md", "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$repoBlobBase/a.md", "", "", "a.md")

# Widen the Status and Latest Handback File columns.
$wsZh.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de detail sheet ----------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-30 16:44:57"
$wsDe.Range("K3").Value = "2016-08-30 16:44:57"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBlobBase/a.md", "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$repoBlobBase/a.md", "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBlobBase/b.md", "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$repoBlobBase/a.md", "", "", "a.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Handback report generated."
